$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Refresh the cached "datetimeFigureOut" date field text from 2/20/2020 to
#    4/3/2020 on every slide master and every slide layout (matches what
#    PowerPoint does to the cached field text on every master/layout when the
#    deck is re-saved on a later date).
# ---------------------------------------------------------------------------
$newDate = "4/3/2020"

for ($di = 1; $di -le $p.Designs.Count; $di++) {
    $design = $p.Designs.Item($di)
    $master = $design.SlideMaster

    # -- slide master's own Date Placeholder -------------------------------
    for ($j = 1; $j -le $master.Shapes.Count; $j++) {
        $sh = $master.Shapes.Item($j)
        if ($sh.Name -like "Date Placeholder*") {
            $sh.TextFrame.TextRange.Text = $newDate
        }
    }

    # -- every custom layout owned by this master ---------------------------
    # NOTE: iterating Shapes directly off $master.CustomLayouts.Item($li) does
    # not reliably persist edits, so visit each layout through a throwaway
    # slide (added then immediately removed) which does persist correctly.
    $layouts = $master.CustomLayouts
    for ($li = 1; $li -le $layouts.Count; $li++) {
        $targetLayout = $layouts.Item($li)
        $tempSlide = $p.Slides.AddSlide($p.Slides.Count + 1, $targetLayout)
        $cl = $tempSlide.CustomLayout
        for ($j = 1; $j -le $cl.Shapes.Count; $j++) {
            $sh = $cl.Shapes.Item($j)
            if ($sh.Name -like "Date Placeholder*") {
                $sh.TextFrame.TextRange.Text = $newDate
            }
        }
        $tempSlide.Delete()
    }
}

# ---------------------------------------------------------------------------
# 2) Remove the red numbered-step ovals (callout bubbles "1".."5") that were
#    overlaid on slides 3-6. (Named-parameter PowerShell functions are not
#    reliable in this host, so the removal loop is inlined per slide.)
# ---------------------------------------------------------------------------

# Slide 3: "Oval 9"
$targetNames = @("Oval 9")
$slide = $p.Slides.Item(3)
for ($j = $slide.Shapes.Count; $j -ge 1; $j--) {
    $sh = $slide.Shapes.Item($j)
    if ($targetNames -contains $sh.Name) {
        $sh.Delete()
    }
}

# Slide 4: "Oval 7"
$targetNames = @("Oval 7")
$slide = $p.Slides.Item(4)
for ($j = $slide.Shapes.Count; $j -ge 1; $j--) {
    $sh = $slide.Shapes.Item($j)
    if ($targetNames -contains $sh.Name) {
        $sh.Delete()
    }
}

# Slide 5: "Oval 6" and "Oval 8"
$targetNames = @("Oval 6", "Oval 8")
$slide = $p.Slides.Item(5)
for ($j = $slide.Shapes.Count; $j -ge 1; $j--) {
    $sh = $slide.Shapes.Item($j)
    if ($targetNames -contains $sh.Name) {
        $sh.Delete()
    }
}

# Slide 6: "Oval 6"
$targetNames = @("Oval 6")
$slide = $p.Slides.Item(6)
for ($j = $slide.Shapes.Count; $j -ge 1; $j--) {
    $sh = $slide.Shapes.Item($j)
    if ($targetNames -contains $sh.Name) {
        $sh.Delete()
    }
}
